$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 178: the script re-ran and found no page number for this term, so the
# previous "NA" placeholder is cleared out.
$ws.Range("C178").Value = ""

# Three new rows of results appended by the script's latest run.
$ws.Range("A179").NumberFormat = "@"
$ws.Range("A179").Value = "2025-07-14"
$ws.Range("A179").Style = "Normal"
$ws.Range("B179").Value = "développement durable"
$ws.Range("C179").Value = 93
$ws.Range("D179").Value = 1

$ws.Range("A180").NumberFormat = "@"
$ws.Range("A180").Value = "2025-07-14"
$ws.Range("A180").Style = "Normal"
$ws.Range("B180").Value = "ruissellement"
$ws.Range("C180").Value = 95
$ws.Range("D180").Value = 2

$ws.Range("A181").NumberFormat = "@"
$ws.Range("A181").Value = "2025-07-14"
$ws.Range("A181").Style = "Normal"
$ws.Range("B181").Value = "eaux de surface"
$ws.Range("C181").Value = 100
$ws.Range("D181").Value = 1
